$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds numeric-looking text (e.g. "43.825.04", "1.00") that must stay
# as literal text, not be reinterpreted as numbers/dates by Excel. Force the
# whole column to Text format for the duration of the update, then restore the
# original (default/General) cell style once all values are written.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '43.825.04'
$ws.Range("E2").Value = '  +4.70%  '
$ws.Range("D3").Value = '2.280.50'
$ws.Range("E3").Value = '  +2.49%  '
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").Value = '231.74'
$ws.Range("E5").Value = '  +0.28%  '
$ws.Range("D6").Value = '0.627'
$ws.Range("E6").Value = '  +0.67%  '
$ws.Range("D7").Value = '61.62'
$ws.Range("E7").Value = '  +1.65%  '
$ws.Range("E8").Value = '  -0.03%  '
$ws.Range("E9").Value = '  +5.13%  '
$ws.Range("B10").Value = 'OKB'
$ws.Range("C10").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D10").Value = '57.78'
$ws.Range("E10").Value = '  -1.75%  '
$ws.Range("B11").Value = 'Dogecoin'
$ws.Range("C11").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D11").Value = '0.0934'
$ws.Range("E11").Value = '  +5.07%  '
$ws.Range("B12").Value = 'TRON'
$ws.Range("C12").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D12").Value = '0.104'
$ws.Range("E12").Value = '  +0.81%  '
$ws.Range("B13").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C13").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D13").Value = '2.618.83'
$ws.Range("E13").Value = '  +2.48%  '
$ws.Range("B14").Value = 'Chainlink'
$ws.Range("C14").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D14").Value = '15.81'
$ws.Range("E14").Value = '  +1.11%  '
$ws.Range("B15").Value = 'Avalanche'
$ws.Range("C15").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D15").Value = '23.90'
$ws.Range("E15").Value = '  +9.94%  '
$ws.Range("B16").Value = 'Polkadot'
$ws.Range("C16").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D16").Value = '5.77'
$ws.Range("E16").Value = '  +4.05%  '
$ws.Range("B17").Value = 'Polygon'
$ws.Range("C17").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D17").Value = '0.813'
$ws.Range("E17").Value = '  +1.98%  '
$ws.Range("B18").Value = 'WrappedEther'
$ws.Range("C18").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D18").Value = '2.287.61'
$ws.Range("E18").Value = '  +3.20%  '
$ws.Range("B19").Value = 'WrappedBTC'
$ws.Range("C19").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D19").Value = '43.664.43'
$ws.Range("E19").Value = '  +4.67%  '
$ws.Range("B20").Value = 'ShibaInu'
$ws.Range("C20").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D20").Value = '0.0₃0935'
$ws.Range("E20").Value = '  +4.85%  '
$ws.Range("B21").Value = 'Litecoin'
$ws.Range("C21").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D21").Value = '73.37'
$ws.Range("E21").Value = '  +0.84%  '
$ws.Range("B22").Value = 'Uniswap'
$ws.Range("C22").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D22").Value = '6.25'
$ws.Range("E22").Value = '  +3.94%  '
$ws.Range("B23").Value = 'BitcoinCash'
$ws.Range("C23").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D23").Value = '250.85'
$ws.Range("E23").Value = '  +0.45%  '
$ws.Range("B24").Value = 'Dai'
$ws.Range("C24").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D24").Value = '1.00'
$ws.Range("E24").Value = '  -0.05%  '
$ws.Range("B25").Value = 'PancakeSwap'
$ws.Range("C25").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D25").Value = '2.55'
$ws.Range("E25").Value = '  +6.86%  '
$ws.Range("B26").Value = 'Toncoin'
$ws.Range("C26").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D26").Value = '2.42'
$ws.Range("E26").Value = '  +4.83%  '
$ws.Range("B27").Value = 'Cosmos'
$ws.Range("C27").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D27").Value = '9.86'
$ws.Range("E27").Value = '  +2.40%  '
$ws.Range("B28").Value = 'Monero'
$ws.Range("C28").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D28").Value = '171.13'
$ws.Range("E28").Value = '  +2.28%  '
$ws.Range("B29").Value = 'Kaspa'
$ws.Range("C29").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D29").Value = '0.142'
$ws.Range("E29").Value = '  -0.26%  '
$ws.Range("B30").Value = 'EthereumClassic'
$ws.Range("C30").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D30").Value = '20.59'
$ws.Range("E30").Value = '  +3.37%  '
$ws.Range("B31").Value = 'ImmutableX'
$ws.Range("C31").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D31").Value = '1.48'
$ws.Range("E31").Value = '  +6.12%  '
$ws.Range("B32").Value = 'WEMIXToken'
$ws.Range("C32").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D32").Value = '2.65'
$ws.Range("E32").Value = '  +0.71%  '
$ws.Range("B33").Value = 'Stellar'
$ws.Range("C33").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D33").Value = '0.122'
$ws.Range("E33").Value = '  +0.56%  '
$ws.Range("D34").Value = '4.79'
$ws.Range("E34").Value = '  +4.09%  '
$ws.Range("B35").Value = 'InternetComputer(DFINITY)'
$ws.Range("C35").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D35").Value = '5.05'
$ws.Range("E35").Value = '  +2.27%  '
$ws.Range("B36").Value = 'Hedera'
$ws.Range("C36").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D36").Value = '0.0662'
$ws.Range("E36").Value = '  +5.99%  '
$ws.Range("B37").Value = 'THORChain'
$ws.Range("C37").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D37").Value = '6.54'
$ws.Range("E37").Value = '  -1.52%  '
$ws.Range("B38").Value = 'LidoDAOToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D38").Value = '2.42'
$ws.Range("E38").Value = '  +2.97%  '
$ws.Range("B39").Value = 'RenderToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D39").Value = '3.65'
$ws.Range("E39").Value = '  -0.72%  '
$ws.Range("B40").Value = 'VeChain'
$ws.Range("C40").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D40").Value = '0.0252'
$ws.Range("E40").Value = '  +5.02%  '
$ws.Range("B41").Value = 'BinanceUSD'
$ws.Range("C41").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D41").Value = '1.00'
$ws.Range("E41").Value = '  -0.20%  '
$ws.Range("B42").Value = 'FraxShare'
$ws.Range("C42").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D42").Value = '8.81'
$ws.Range("E42").Value = '  +2.50%  '
$ws.Range("B43").Value = 'TerraClassic'
$ws.Range("C43").Value = 'https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc'
$ws.Range("D43").Value = '0.000224'
$ws.Range("E43").Value = '  -12.24%  '
$ws.Range("D44").Value = '1.23'
$ws.Range("E44").Value = '  +1.01%  '
$ws.Range("B45").Value = 'FTXToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D45").Value = '4.50'
$ws.Range("E45").Value = '  -6.14%  '
$ws.Range("B46").Value = 'Cronos'
$ws.Range("C46").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D46").Value = '0.0969'
$ws.Range("E46").Value = '  -1.08%  '
$ws.Range("B47").Value = 'Aave'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D47").Value = '98.43'
$ws.Range("E47").Value = '  -0.18%  '
$ws.Range("B48").Value = 'Maker'
$ws.Range("C48").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D48").Value = '1.474.27'
$ws.Range("E48").Value = '  +0.43%  '
$ws.Range("B49").Value = 'InjectiveProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D49").Value = '16.74'
$ws.Range("E49").Value = '  +1.66%  '
$ws.Range("B50").Value = 'NEARProtocol'
$ws.Range("C50").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D50").Value = '2.30'
$ws.Range("E50").Value = '  +10.42%  '
$ws.Range("E51").Value = '  +1.69%  '

# Restore default (General) style/number format on column D
$ws.Range("D2:D51").Style = "Normal"
